# Header updates for summer uploads
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (row 1) to the new header labels.
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("G1").Value = "Ministry Course Code and Level"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"

# Update the selection shown in the saved worksheet view.
$ws.Range("A1:K1").Select()
